$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 36: new "ROM 1 / IP Core" entry (3.4.2020, 13:30-14:00)
# ---------------------------------------------------------------------
$a36 = $ws.Range("A36")
$a36.NumberFormat = "@"
$a36.HorizontalAlignment = -4152
$a36.VerticalAlignment = -4108
$a36.Value = "3.4.2020"

$b36 = $ws.Range("B36")
$b36.NumberFormat = "h:mm"
$b36.Value = 0.5625

$c36 = $ws.Range("C36")
$c36.NumberFormat = "h:mm"
$c36.Value = 0.58333333333333337

$d36 = $ws.Range("D36")
$d36.NumberFormat = "h:mm"
$d36.HorizontalAlignment = -4152
$d36.VerticalAlignment = -4108
$d36.Formula = "=C36-B36"

$e36 = $ws.Range("E36")
$e36.HorizontalAlignment = -4131
$e36.VerticalAlignment = -4108
$e36.Value = "ROM 1"

$f36 = $ws.Range("F36")
$f36.HorizontalAlignment = -4131
$f36.VerticalAlignment = -4108
$f36.Value = "IP Core"

# ---------------------------------------------------------------------
# Row 37: new "Mem Ctrl 1 / Arch and TB" entry (3.4.2020, 14:00-15:00)
# ---------------------------------------------------------------------
$a37 = $ws.Range("A37")
$a37.NumberFormat = "@"
$a37.HorizontalAlignment = -4152
$a37.VerticalAlignment = -4108
$a37.Value = "3.4.2020"

$b37 = $ws.Range("B37")
$b37.NumberFormat = "h:mm"
$b37.Value = 0.58333333333333337

$c37 = $ws.Range("C37")
$c37.NumberFormat = "h:mm"
$c37.Value = 0.625

$d37 = $ws.Range("D37")
$d37.NumberFormat = "h:mm"
$d37.Formula = "=C37-B37"

$e37 = $ws.Range("E37")
$e37.HorizontalAlignment = -4131
$e37.VerticalAlignment = -4108
$e37.Value = "Mem Ctrl 1"

$f37 = $ws.Range("F37")
$f37.HorizontalAlignment = -4131
$f37.VerticalAlignment = -4108
$f37.Value = "Arch and TB"

# ---------------------------------------------------------------------
# Row 38: new "VGA Top / Arch and TB" entry (4.4.2020, starting 10:00)
# ---------------------------------------------------------------------
$a38 = $ws.Range("A38")
$a38.NumberFormat = "@"
$a38.HorizontalAlignment = -4152
$a38.VerticalAlignment = -4108
$a38.Value = "4.4.2020"

$b38 = $ws.Range("B38")
$b38.NumberFormat = "h:mm"
$b38.Value = 0.41666666666666669

$e38 = $ws.Range("E38")
$e38.HorizontalAlignment = -4131
$e38.VerticalAlignment = -4108
$e38.Value = "VGA Top"

$f38 = $ws.Range("F38")
$f38.HorizontalAlignment = -4131
$f38.VerticalAlignment = -4108
$f38.Value = "Arch and TB"

$ws.Range("G38").Value = "Add VGA PLL"

# ---------------------------------------------------------------------
# Update the view: scroll so row 13 is at the top, select A39
# ---------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$ws.Range("A39").Select()
